# Mifos "ProductLoanInput" loan-product workbook:
#  - Row 6 used to hold the generic "Currency" / "US Dollar " pair; it is
#    re-labelled to the lowercase "currency" / "US Dollar" (no trailing
#    space) strings used elsewhere in the automation-script vocabulary, and
#    the answer cell (B6) is re-painted with the same green "answered"
#    highlight used on other verified rows.
#  - The active sheet/selection moves from ProductLoanOutput!B13 to
#    ProductLoanInput!A6:B6 (the field that was just edited), scrolling the
#    input sheet back to the top.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")

# Relabel the currency question/answer pair.
$ws1.Range("A6").Value = "currency"

$answer = $ws1.Range("B6")
$answer.ClearFormats()
$answer.Interior.Color = 5296274   # BGR for RGB FF92D050 (the "answered" green fill)
$answer.Value = "US Dollar"

# Move the active sheet/selection to the input sheet, scrolled to the top,
# with the currency row selected.
$ws1.Activate()
$ws1.Range("A6:B6").Select()
